# Auto-generated Excel COM-interop script to apply cryptos.xlsx price/volume update
# (commit: "Updated cryptos list on Sat Aug 12 13:42:47 UTC 2023 with GitHub Actions")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (non-numeric) formatting for touched cells so that values such as
# "1.000", "29.407.73" or "0.00000000119" are stored verbatim as strings rather
# than being auto-converted by Excel into numbers/scientific notation.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.407.73'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.851.38'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.24%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.87'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6289'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.41%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07694'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.83%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2936'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.57'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07749'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.75%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.849.67'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.00001108'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +10.25%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6815'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.45%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.67'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.77%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.106.48'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.151'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '29.450.04'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '229.01'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.435'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.15'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1386'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.42%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.388'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.17%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.69'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.315'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +4.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.468'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05705'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.43%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.053'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.73%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.848'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.12%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7087'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.38%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.586'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.81%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.219.47'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.487'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +5.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9108'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.0000'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.015.16'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.70'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '66.37'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.45%  '
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.131'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.60%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000119'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.47%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4016'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.014'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.83%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.681'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.04%  '
